$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the capacity column (D2:D7) from 10000 to 250
$ws.Range("D2:D7").Value = 250
